$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.046.60"
$ws.Range("E2").Value = "  +1.36%  "

$ws.Range("D3").Value = "'2.245.13"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'318.24"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").Value = "'101.08"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  -1.72%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.543"
$ws.Range("E9").Value = "  -3.71%  "

$ws.Range("D10").Value = "'36.74"
$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("D14").Value = "'2.588.10"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "'2.280.92"
$ws.Range("E15").Value = "  +2.85%  "

$ws.Range("D16").Value = "'0.848"
$ws.Range("E16").Value = "  -2.54%  "

$ws.Range("D17").Value = "'14.16"
$ws.Range("E17").Value = "  -1.84%  "

$ws.Range("D18").Value = "'43.911.71"
$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("D19").Value = "'13.42"
$ws.Range("E19").Value = "  -5.39%  "

$ws.Range("D20").Value = "'0.0₃0976"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").Value = "'6.43"
$ws.Range("E21").Value = "  -3.10%  "

$ws.Range("D22").Value = "'65.38"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("E23").Value = "  -4.56%  "

$ws.Range("D24").Value = "'234.98"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = "  -6.32%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "  +3.67%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'37.48"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("D30").Value = "'6.05"
$ws.Range("E30").Value = "  -5.79%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.08"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'158.22"
$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("D33").Value = "'0.0846"
$ws.Range("E33").Value = "  -4.16%  "

$ws.Range("D34").Value = "'2.68"
$ws.Range("E34").Value = "  -1.80%  "

$ws.Range("D35").Value = "'3.19"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  +6.98%  "

$ws.Range("D37").Value = "'1.94"
$ws.Range("E37").Value = "  +2.77%  "

$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("D39").Value = "'15.87"
$ws.Range("E39").Value = "  +9.71%  "

$ws.Range("D40").Value = "'3.67"
$ws.Range("E40").Value = "  -2.26%  "

$ws.Range("D41").Value = "'4.13"
$ws.Range("E41").Value = "  -7.09%  "

$ws.Range("E42").Value = "  -3.21%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "'1.750.37"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.196"
$ws.Range("E45").Value = "  -4.15%  "

$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").Value = "'74.49"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "'81.39"
$ws.Range("E47").Value = "  -4.28%  "

$ws.Range("D48").Value = "'5.14"
$ws.Range("E48").Value = "  -3.27%  "

$ws.Range("D49").Value = "'102.68"
$ws.Range("E49").Value = "  -1.29%  "

$ws.Range("D50").Value = "'1.66"
$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("D51").Value = "'57.35"
$ws.Range("E51").Value = "  -2.78%  "
